$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (pushes existing rows 2-13 down to 3-14)
$ws.Rows.Item(2).Insert()

# Copy formatting/values from the row that is now row 3 (the old "LeftSide" row)
# into the newly inserted row 2, so the new row matches the existing row styling.
$ws.Range("A3:L3").Copy($ws.Range("A2:L2"))

# Set the new row's field name to "Tag"
$ws.Range("A2").Value = "Tag"

# Update the saved selection to match the target workbook state
$ws.Range("A3").Select()

Write-Host "done"
